# Daily update at 8 AM UTC
# Adds the next day's row (row 72) to the "Wins Over Time" tracking sheet,
# and moves the "last row" date-only formatting from A71 to the new A72.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Previously the last row (71) used the date-only number format to mark it
# as the final entry. Since row 72 becomes the new final entry, row 71's
# date cell goes back to the regular timestamp format used by all the
# other (non-last) rows.
$ws.Range("A71").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New data row for the next day.
$ws.Range("A72").Value = 45812
$ws.Range("B72").Value = 307
$ws.Range("C72").Value = 304
$ws.Range("D72").Value = 309

# The newest row's date cell gets the "last row" date-only number format.
$ws.Range("A72").NumberFormat = "YYYY-MM-DD"
